$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-01-28 Sunday" "2024-01-29 Monday"

Replace-Text "50×14=" "34×30="
Replace-Text "59×55=" "99×30="
Replace-Text "50×58=" "77×98="
Replace-Text "84×62=" "46×31="
Replace-Text "55×11=" "23×87="

Replace-Text "86×14=" "36×98="
Replace-Text "85×91=" "11×20="
Replace-Text "73×29=" "67×51="
Replace-Text "68×38=" "76×37="
Replace-Text "98×96=" "36×59="

Replace-Text "90×72=" "51×78="
Replace-Text "97×29=" "48×70="
Replace-Text "46×93=" "83×49="
Replace-Text "50×63=" "96×16="
Replace-Text "82×48=" "23×11="

Replace-Text "84×41=" "56×30="
Replace-Text "28×45=" "22×86="
Replace-Text "64×92=" "47×30="
Replace-Text "18×50=" "58×84="
Replace-Text "82×97=" "58×99="

Replace-Text "29×14=" "87×61="
Replace-Text "88×44=" "77×73="
Replace-Text "71×19=" "61×55="
Replace-Text "66×34=" "76×47="
Replace-Text "76×86=" "28×96="

Write-Output "Done"
